# Leave card update (1/9/2024 4:42 pm):
#   - Insert a new leave-card row (row 45) on Sheet1 recording a
#     "UT(0-0-24)" undertime entry of 0.05 (5%) in the
#     "Absence Undertime W/ Pay" column, pushing all subsequent rows down
#     by one (old row 45 -> 46, ..., old row 129 -> 130).
#   - Grow Table1 to match the extra row.
#   - CONVERTION!F3 (days-undertime lookup key) changes from 35 to 24,
#     which ripples into CONVERTION!G3 and Sheet1!E9 via formulas.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Grow the table first so the shifted-down last row keeps its
# "totals-row-ish" bottom border style (dxf ids 11-14) instead of turning
# into an ordinary interior row.
$lo.Resize($ws.Range("A8:K130"))

# Insert a new blank row above the old row 45 (1/7/2023 entry), shifting
# rows 45:129 down to 46:130.
$ws.Rows("45:45").Insert()

# New row 45 should look like a normal interior data row: pick up the
# same cell formatting (borders/number formats) the row below already has.
$ws.Range("A46:K46").Copy()
$ws.Range("A45:K45").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new leave-card entry.
$ws.Range("B45").Value = "UT(0-0-24)"
$ws.Range("D45").Value = 0.05
$ws.Range("G45").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Conversion table lookup key used for the undertime balance formula.
$wsConv = $wb.Worksheets.Item("CONVERTION")
$wsConv.Range("F3").Value = 24

# Match the author's final cursor position on Sheet1.
$ws.Select()
$ws.Range("F49").Select()

$wb.Application.Calculate() | Out-Null
